# Add a new row of data (email + password) below the existing header/data
# row on Sheet1, turning the email cell into a mailto: hyperlink (Excel's
# built-in "Hyperlink" cell style - underlined, theme colour 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$email = "amzftv@gmail.com"
$password = "ab410410"

# New data row
$ws.Range("A2").Value = $email
$ws.Range("B2").Value = $password

# Turn A2 into a live mailto: hyperlink (this also applies the built-in
# "Hyperlink" cell style: underlined font, theme color 10).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:" + $email)

# Match the saved selection (B2 active/selected, as recorded in the sheet view).
$ws.Range("B2").Select() | Out-Null
